$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 'Karamaji Eliphasi', 7, 18, 1, 3, 16),
    @(3, 'Katimbo Dennis Ronald', 3, 28, 0, 5, 2),
    @(4, 'Ekiru David Bernard', 6, 17, 0, 15, 4),
    @(5, 'Muwanguzi Enock', 1, 12, 0, 0, 0),
    @(6, 'Ssali Musa', 5, 14, 0, 2, 1),
    @(7, 'Kizza Harper Stephen', 2, 30, 2, 1, 1),
    @(8, 'Kwesiga Norman', 0, 1, 1, 1, 0),
    @(9, 'Ssekanyonyi Kennedy', 8, 37, 3, 0, 2),
    @(10, 'kibalama frank', 13, 49, 2, 2, 16),
    @(11, 'Mukasa George Parker', 4, 20, 0, 6, 7),
    @(12, 'Barya Mwebaze', 7, 40, 0, 4, 13),
    @(13, 'Makumbi Patrick', 6, 25, 3, 2, 0),
    @(14, 'Sserunyigo Enock', 14, 37, 2, 1, 5),
    @(15, 'David Oluka', 1, 21, 0, 0, 20),
    @(16, 'Kiddawalime Jimmy', 0, 1, 0, 0, 0),
    @(17, 'Mukasa joseph mutesasira', 3, 13, 0, 0, 0),
    @(18, 'Wateya Job', 5, 21, 0, 1, 4),
    @(19, 'Musaazi Brian', 7, 6, 0, 3, 1),
    @(20, 'Mubiru Stephen', 2, 6, 0, 0, 0),
    @(21, 'Kwesiga Samuel', 6, 20, 0, 17, 0),
    @(22, 'Aojan James Patrick', 0, 0, 0, 0, 0),
    @(23, 'Kyeyune Eric Brian', 2, 10, 0, 2, 6),
    @(24, 'Musasizi Eric John', 7, 22, 0, 1, 1),
    @(25, 'Muyanja Moses', 5, 38, 0, 0, 0),
    @(26, 'MBAZIIRA HUDHAIFA', 12, 29, 0, 0, 0),
    @(27, 'KIWANDA PETER', 0, 2, 0, 0, 3),
    @(28, 'Daphine Nakanwagi', 4, 26, 0, 0, 0),
    @(29, 'Lubega Mark', 2, 4, 0, 0, 0),
    @(30, 'Talemwa Steven', 0, 11, 0, 0, 0),
    @(31, 'other', 0, 0, 0, 0, 2),
    @(32, 'Bbaale  Jean', 4, 29, 0, 2, 0),
    @(33, 'Boas Massinde', 0, 0, 0, 0, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}

# Row 34 (Boas Massinde) no longer exists in the report; delete the row entirely
$ws.Rows.Item(34).Delete()

